$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Text change: "Ready for handoff" -> "In Translation"
#    This string is used by the Status / status-like columns on all three
#    sheets: Overview (E2:F4), zh-cn (C2:C4) and de-de (C2:C4).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = "In Translation"

# ---------------------------------------------------------------------------
# 2. Column width changes (report regenerated with narrower Status columns):
#    - Overview: columns E and F (5, 6)
#    - zh-cn:    column C (3)
#    - de-de:    column C (3)
#    Original width 17.2159881591797 -> new width 13.4101845877511
# ---------------------------------------------------------------------------
$newWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newWidth
